$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("D2").Value = 91.38
$ws.Range("F2").Value = 3.84
$ws.Range("H2").Value = 73
$ws.Range("N2").Value = 85.77505782882612

# Row 3 update
$ws.Range("N3").Value = 85.77505782882612
